# Added 2026 Inductee Photos
# Adds three new inductees (rows 149-151) to the "All Inductees" sheet:
#   Fred Zuercher, George Demetriou, Terry Angell - all class of 2026

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 149 - Fred Zuercher
$ws.Range("A149").Value = "Fred Zuercher"
$ws.Range("B149").Value = 2026
$ws.Range("C149").Value = "Fred_Zuercher.jpg"
$ws.Range("D149").Value = "Fred_Zuercher"

# Row 150 - George Demetriou
$ws.Range("A150").Value = "George Demetriou"
$ws.Range("B150").Value = 2026
$ws.Range("C150").Value = "George_Demetriou.jpg"
$ws.Range("D150").Value = "George_Demetriou"

# Row 151 - Terry Angell
$ws.Range("A151").Value = "Terry Angell"
$ws.Range("B151").Value = 2026
$ws.Range("C151").Value = "Terry_Angell.jpg"
$ws.Range("D151").Value = "Terry_Angell"

# Scroll the view down to the newly-added rows and move the active selection
# to C154 (matching the author's on-save cursor position).
$win = $excel.ActiveWindow
$win.ScrollRow = 133
$win.ScrollColumn = 1
$ws.Range("C154").Select()
